$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Terminal La Palmera de La Serena - Espinaca".
# It is inserted as a new data row right above the current row 372, which pushes all
# subsequent rows (old 372..421) down by one (new 373..422), growing the sheet from
# A1:R421 to A1:R422.
$ws.Rows.Item(372).Insert()

# Populate the newly inserted row 372 with the new weekly record.
$ws.Cells.Item(372, 1).Value2 = 8
$ws.Cells.Item(372, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(372, 3).Value2 = "Coquimbo"
$ws.Cells.Item(372, 4).Value2 = 45077
$ws.Cells.Item(372, 5).Value2 = 4
$ws.Cells.Item(372, 6).Value2 = 100112012
$ws.Cells.Item(372, 7).Value2 = "Espinaca"
$ws.Cells.Item(372, 8).Value2 = "Sin especificar"
$ws.Cells.Item(372, 9).Value2 = "Primera"
$ws.Cells.Item(372, 10).Value2 = 1200
$ws.Cells.Item(372, 11).Value2 = 400
$ws.Cells.Item(372, 12).Value2 = 500
$ws.Cells.Item(372, 13).Value2 = 450
$ws.Cells.Item(372, 14).Value2 = "$/atado 300 a 500 gramos"
$ws.Cells.Item(372, 15).Value2 = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(372, 16).Value2 = 900
$ws.Cells.Item(372, 17).Value2 = 0.5
$ws.Cells.Item(372, 18).Value2 = "Hortaliza"
